$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.956.41"
$ws.Range("E2").Value = "  +2.06%  "

$ws.Range("D3").Value = "2.289.16"
$ws.Range("E3").Value = "  +3.13%  "

$ws.Range("D5").Value = "'252.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("E6").Value = "  +3.85%  "

$ws.Range("D7").Value = "'73.68"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.98%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.645"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.93%  "

$ws.Range("D10").Value = "'39.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("D11").Value = "'0.0973"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.44%  "

$ws.Range("D12").Value = "'59.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.58%  "

$ws.Range("D13").Value = "'7.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.98%  "

$ws.Range("E14").Value = "  +1.57%  "

$ws.Range("D15").Value = "2.634.37"
$ws.Range("E15").Value = "  +3.13%  "

$ws.Range("D16").Value = "'15.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.42%  "

$ws.Range("D17").Value = "'0.873"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("D18").Value = "2.299.72"
$ws.Range("E18").Value = "  +3.26%  "

$ws.Range("D19").Value = "42.871.98"
$ws.Range("E19").Value = "  +1.99%  "

$ws.Range("D20").Value = "'0.0000100"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.47%  "

$ws.Range("D21").Value = "'6.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.87%  "

$ws.Range("D22").Value = "'72.46"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "'237.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.01%  "

$ws.Range("E24").Value = "  +9.20%  "

$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("D26").Value = "'11.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.02%  "

$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("E28").Value = "  +0.60%  "

$ws.Range("D29").Value = "'3.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("D30").Value = "'2.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.16%  "

$ws.Range("D31").Value = "'167.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.09%  "

$ws.Range("D32").Value = "'21.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.95%  "

$ws.Range("E33").Value = "  +5.01%  "

$ws.Range("E34").Value = "  +5.61%  "

$ws.Range("D35").Value = "'0.0823"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.35%  "

$ws.Range("D36").Value = "'31.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.31%  "

$ws.Range("E37").Value = "  +3.43%  "

$ws.Range("D38").Value = "'4.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.51%  "

$ws.Range("D39").Value = "'4.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.23%  "

$ws.Range("D40").Value = "'0.0309"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.09%  "

$ws.Range("D41").Value = "'14.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +17.14%  "

$ws.Range("E42").Value = "  +4.49%  "

$ws.Range("D43").Value = "'5.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.26%  "

$ws.Range("D44").Value = "'0.218"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.13%  "

$ws.Range("D45").Value = "'61.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("E46").Value = "  +5.50%  "

$ws.Range("D47").Value = "'4.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.18%  "

$ws.Range("E48").Value = "  +3.46%  "

$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("E50").Value = "  +1.76%  "

$ws.Range("E51").Value = "  +6.60%  "
